$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update weights for recent player performance
$ws.Range("C5").Value = 18.9
$ws.Range("M6").Value = 12.5

# Mark the stacks that hit "Success" once row8 value formulas produce data
$ws.Range("D8").Value = "Success"
$ws.Range("N8").Value = "Success"

# Update the active selection cursor
$ws.Range("N9").Select()
